$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 588 ("人生は美しい、そう望んだ者には。" post),
# which shifts all subsequent rows up by one.
$ws.Rows.Item(588).Delete()
